$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (text moved up rows in B/C columns; "status" relocated to column A) ---

# Row 3: B3 and C3 change from "lozinka" to "telefon"
$ws.Range("B3").Value = "telefon"
$ws.Range("C3").Value = "telefon"

# Row 4: B4 and C4 change from "telefon" to "mejl"
$ws.Range("B4").Value = "mejl"
$ws.Range("C4").Value = "mejl"

# Row 5: A5 gains "status"; B5/C5 become "slika"
$ws.Range("A5").Value = "status"
$ws.Range("B5").Value = "slika"
$ws.Range("C5").Value = "slika"

# Row 6: B6 -> "ime"; C6 -> "naziv"
$ws.Range("B6").Value = "ime"
$ws.Range("C6").Value = "naziv"

# Row 7: B7 -> "prezime"; C7 -> "ulica"
$ws.Range("B7").Value = "prezime"
$ws.Range("C7").Value = "ulica"

# Row 8: B8 cleared; C8 -> "grad"
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = "grad"

# Row 9: B9 cleared; C9 -> "drzava"
$ws.Range("B9").ClearContents()
$ws.Range("C9").Value = "drzava"

# Row 10: C10 -> "maticniBroj"
$ws.Range("C10").Value = "maticniBroj"

# Row 11: C11 -> "opis"
$ws.Range("C11").Value = "opis"

# Row 12: C12 cleared
$ws.Range("C12").ClearContents()

# Row 13: C13 cleared
$ws.Range("C13").ClearContents()

# --- Restyle column A (rows 2-11) from style index 6 to style index 1 ---
# Style 6 = fontId0 fillId0 borderId1 applyFill applyBorder + center alignment (duplicate of 1 minus applyFill)
# Style 1 = fontId0 fillId0 borderId1 applyBorder + center alignment
# Re-applying the (already present) thin border re-derives the cell format and
# collapses it onto the equivalent existing style record (index 1).
$ws.Range("A2:A11").Borders.LineStyle = 1

# --- Selection change ---
$ws.Range("C12:C13").Select()
